$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D4: 100 -> 2
$ws.Range("D4").Value = 2

# D5: add value 3
$ws.Range("D5").Value = 3

# C6: new shared string "answer answer"
$ws.Range("C6").Value = "answer answer"

# Update selection to C6
$ws.Range("C6").Select()
